$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Content: write the new text into A1
$ws.Range("A1").Value = "First change"

# Formatting: Calibri font (explicit) + vertically centered alignment
$ws.Range("A1").Font.Name = "Calibri"
$ws.Range("A1").VerticalAlignment = -4108

# Page setup: portrait, paper size 9 (A4)
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9
